# Auto-generated edit script applying scheduled market-data refresh to Phantom_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5053.875
$ws.Range("I64").Value = 4826.2
$ws.Range("J64").Value = 5433.3335
$ws.Range("K64").Value = 4826.2
$ws.Range("L64").Value = 5433.3335
$ws.Range("M64").Value = -4578.2
$ws.Range("N64").Value = -5929.3335
$ws.Range("H67").Value = 5053.875
$ws.Range("I67").Value = 4826.2
$ws.Range("J67").Value = 5433.3335
$ws.Range("K67").Value = 4826.2
$ws.Range("L67").Value = 5433.3335
$ws.Range("M67").Value = -3968.2
$ws.Range("N67").Value = -7149.3335
$ws.Range("H92").Value = 629.6923
$ws.Range("I92").Value = 686.625
$ws.Range("K92").Value = 686.625
$ws.Range("M92").Value = 561.375
$ws.Range("H94").Value = 450
$ws.Range("I94").Value = 450
$ws.Range("K94").Value = 450
$ws.Range("M94").Value = 1
$ws.Range("H107").Value = 1246.15
$ws.Range("I107").Value = 1206.4736
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1206.4736
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 713.5264
$ws.Range("N107").Value = -5840
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H116").Value = 4803.6665
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558
$ws.Range("H132").Value = 5340.905
$ws.Range("I132").Value = 5703.706
$ws.Range("J132").Value = 3799
$ws.Range("K132").Value = 17111.118
$ws.Range("L132").Value = 11397
$ws.Range("M132").Value = -14581.118
$ws.Range("N132").Value = -16457
$ws.Range("H138").Value = 2481.1765
$ws.Range("I138").Value = 2373
$ws.Range("J138").Value = 2514.4614
$ws.Range("K138").Value = 7119
$ws.Range("L138").Value = 7543.3842
$ws.Range("M138").Value = -1979
$ws.Range("N138").Value = -17823.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 4249.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H102").Value = 1224.75
$ws.Range("I102").Value = 1224.75
$ws.Range("K102").Value = 1224.75
$ws.Range("M102").Value = 397.25
$ws.Range("H104").Value = 25460.25
$ws.Range("J104").Value = 25460.25
$ws.Range("L104").Value = 25460.25
$ws.Range("N104").Value = -32448.25
$ws.Range("H110").Value = 5156.5454
$ws.Range("I110").Value = 5934.4443
$ws.Range("J110").Value = 1656
$ws.Range("K110").Value = 5934.4443
$ws.Range("L110").Value = 1656
$ws.Range("M110").Value = -3889.4443
$ws.Range("N110").Value = -5746
$ws.Range("H122").Value = 1889.1111
$ws.Range("I122").Value = 1889.1111
$ws.Range("K122").Value = 5667.3333
$ws.Range("M122").Value = -3217.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2677.75
$ws.Range("I107").Value = 2677.75
$ws.Range("K107").Value = 2677.75
$ws.Range("M107").Value = -757.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H32").Value = 8607
$ws.Range("I32").Value = 2955
$ws.Range("J32").Value = 19911
$ws.Range("K32").Value = 2955
$ws.Range("L32").Value = 19911
$ws.Range("M32").Value = -2639
$ws.Range("N32").Value = -20543

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 214.88889
$ws.Range("J2").Value = 490.14285
$ws.Range("L2").Value = 2940.8571
$ws.Range("N2").Value = -3166.8571
$ws.Range("H44").Value = 63
$ws.Range("I44").Value = 63
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 189
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 209
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 17010
$ws.Range("I52").Value = 17010
$ws.Range("K52").Value = 17010
$ws.Range("M52").Value = -16751

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 2568
$ws.Range("I11").Value = 200
$ws.Range("J11").Value = 3041.6
$ws.Range("K11").Value = 200
$ws.Range("L11").Value = 3041.6
$ws.Range("M11").Value = -60
$ws.Range("N11").Value = -3321.6
$ws.Range("H13").Value = 975
$ws.Range("I13").Value = 126
$ws.Range("J13").Value = 1399.5
$ws.Range("K13").Value = 126
$ws.Range("L13").Value = 1399.5
$ws.Range("M13").Value = 14
$ws.Range("N13").Value = -1679.5
$ws.Range("H19").Value = 4832.6665
$ws.Range("I19").Value = 999
$ws.Range("J19").Value = 6749.5
$ws.Range("K19").Value = 999
$ws.Range("L19").Value = 6749.5
$ws.Range("M19").Value = -829
$ws.Range("N19").Value = -7089.5
$ws.Range("H48").Value = 40000
$ws.Range("I48").Value = 40000
$ws.Range("K48").Value = 40000
$ws.Range("M48").Value = -39339
$ws.Range("H61").Value = 2923.5
$ws.Range("I61").Value = 2414.25
$ws.Range("K61").Value = 2414.25
$ws.Range("M61").Value = -2212.25
$ws.Range("H113").Value = 2923.5
$ws.Range("I113").Value = 2414.25
$ws.Range("K113").Value = 2414.25
$ws.Range("M113").Value = -244.25
$ws.Range("H132").Value = 1912.3334
$ws.Range("I132").Value = 1878.4286
$ws.Range("J132").Value = 2031
$ws.Range("K132").Value = 5635.2858
$ws.Range("L132").Value = 6093
$ws.Range("M132").Value = -3105.2858
$ws.Range("N132").Value = -11153
$ws.Range("H136").Value = 2797
$ws.Range("I136").Value = 2974.2222
$ws.Range("K136").Value = 8922.6666
$ws.Range("M136").Value = -6372.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 29748.25
$ws.Range("I2").Value = 29748.25
$ws.Range("K2").Value = 29748.25
$ws.Range("M2").Value = -29636.25
$ws.Range("H100").Value = 90910820
$ws.Range("I100").Value = 111112610
$ws.Range("J100").Value = 2750
$ws.Range("K100").Value = 222225220
$ws.Range("L100").Value = 5500
$ws.Range("M100").Value = -222224679
$ws.Range("N100").Value = -6582
$ws.Range("H113").Value = 922.875
$ws.Range("I113").Value = 1063.8334
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 3191.5002
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -1021.5002
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 1480.4
$ws.Range("I132").Value = 1475.5
$ws.Range("K132").Value = 4426.5
$ws.Range("M132").Value = -1896.5

Write-Host "Applied all Phantom_Profits cell updates."